$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing revised data points ---
$ws.Cells.Item(434, 2).Value = 1182578000000
$ws.Cells.Item(434, 4).Value = 132326783636.2009

$ws.Cells.Item(435, 2).Value = 1180145000000
$ws.Cells.Item(435, 4).Value = 132200988024.9583

$ws.Cells.Item(451, 2).Value = 1297962000000
$ws.Cells.Item(451, 4).Value = 129240465996.2163

$ws.Cells.Item(457, 2).Value = 1375618000000
$ws.Cells.Item(457, 4).Value = 132124861931.518

$ws.Cells.Item(458, 2).Value = 1365405000000
$ws.Cells.Item(458, 4).Value = 134061699182.1225

# --- Append new rows 460:462, carrying the same formatting as row 459 ---
$ws.Range("A459:D459").Copy($ws.Range("A460:D460"))
$ws.Range("A459:D459").Copy($ws.Range("A461:D461"))
$ws.Range("A459:D459").Copy($ws.Range("A462:D462"))

$ws.Cells.Item(460, 1).Value = 44986
$ws.Cells.Item(460, 2).Value = 1393379000000
$ws.Cells.Item(460, 3).Value = 0.09767818943708059
$ws.Cells.Item(460, 4).Value = 136102737919.6499

$ws.Cells.Item(461, 1).Value = 45017
$ws.Cells.Item(461, 2).Value = 1391394000000
$ws.Cells.Item(461, 3).Value = 0.09943323058566172
$ws.Cells.Item(461, 4).Value = 138350800437.5062

$ws.Cells.Item(462, 1).Value = 45047
$ws.Cells.Item(462, 2).Value = 1394414000000
$ws.Cells.Item(462, 3).Value = 0.09827237170541873
$ws.Cells.Item(462, 4).Value = 137032370919.2398
